# Apply regimen data update & bug fix to the temple.xlsx workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column D ("计算量") before the current "最小值" column.
# This shifts old D (最小值), E (最大值), F (单位) to E, F, G respectively.
$ws.Columns.Item(4).Insert()

# Update header row
$ws.Range("C1").Value = "推荐值(mg/m2)"
$ws.Range("D1").Value = "计算量"
$ws.Range("E1").Value = "最小值"
$ws.Range("F1").Value = "最大值"
$ws.Range("G1").Value = "单位"

# Fill the new "计算量" column with 0 for each data row
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0

# Update selected cell / active cell to B9 (as recorded after the edit)
$ws.Range("B9").Select()

# Resize the window (reflects the author's local window geometry change)
$excel.ActiveWindow.Width = 9615
$excel.ActiveWindow.Height = 7815
